$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41: Max subset sum tree
$ws.Range("A41").Value = 44
$ws.Range("B41").Value = "Max subset sum tree"
$ws.Range("H41").Value = "MaxSubsetSumTree"

# Row 42: NodesAtDistanceK for binary tree
$ws.Range("A42").Value = 45
$ws.Range("B42").Value = "NodesAtDistanceK for binary tree"
$ws.Range("H42").Value = "NodesAtDistanceK"

# Apply special font style to H42 (JetBrains Mono, color FFA9B7C6, size 9.8, vertical center)
$h42 = $ws.Range("H42")
$h42.Font.Name = "JetBrains Mono"
$h42.Font.Size = 9.8
$h42.Font.Color = 11706793
$h42.VerticalAlignment = -4108

# Update selection/view to match final state
$ws.Range("H42").Select()
$ws.Application.ActiveWindow.ScrollRow = 22
